# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-22, columns E/F) gets re-ordered so the
# periods run ascending (2102 -> 2108) instead of descending (2108 -> 2102).
# Only the data in columns E (period label) and F (period end date serial)
# move; every cell keeps its own formatting (the thicker border that marks
# the last data row stays on row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E - "Periodo Mora" labels, now ascending top to bottom.
$ws.Range("E16").Value = "2102"
$ws.Range("E17").Value = "2103"
$ws.Range("E18").Value = "2104"
$ws.Range("E19").Value = "2105"
$ws.Range("E20").Value = "2106"
$ws.Range("E21").Value = "2107"
$ws.Range("E22").Value = "2108"

# Column F - matching period-end date serials; only the first/last rows
# actually swap value (31495 <-> 36341), the rest already held 36341.
$ws.Range("F16").Value = 36341
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 36341
$ws.Range("F22").Value = 31495
